$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "groups" column (E) with the re-classified values ---
# (species were re-grouped as part of the model-averaging update; the
#  "amphib" label is retired and a new "macrophyte" label is introduced)
$ws.Range("E3").Value  = "fish"
$ws.Range("E4").Value  = "macrophyte"
$ws.Range("E5").Value  = "fish"
$ws.Range("E6").Value  = "invert"
$ws.Range("E8").Value  = "algae"
$ws.Range("E9").Value  = "fish"
$ws.Range("E10").Value = "fish"
$ws.Range("E12").Value = "invert"
$ws.Range("E13").Value = "invert"
$ws.Range("E15").Value = "algae"
$ws.Range("E16").Value = "invert"
$ws.Range("E17").Value = "invert"
$ws.Range("E18").Value = "algae"
$ws.Range("E19").Value = "algae"
$ws.Range("E20").Value = "macrophyte"
$ws.Range("E21").Value = "algae"

# --- Page setup: portrait orientation (adds <pageSetup .../>) ---
$ws.PageSetup.Orientation = 1

# --- Move/update the active selection to H10 ---
$ws.Range("H10").Select()
